# Update to 2025 MLB Constants
#  - wOBA & FIP Constants: add AL PA / NL PA columns (Q, R) for the 2025 row
#  - Park Factors: insert a new "League" column after "Abbreviation"
#    and populate it with each team's league (AL/NL)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "wOBA & FIP Constants" sheet: new AL PA / NL PA columns
# ---------------------------------------------------------------------------
$wobaSheet = $wb.Worksheets.Item("wOBA & FIP Constants")

$wobaSheet.Range("Q1").Value = "AL PA"
$wobaSheet.Range("R1").Value = "NL PA"
$wobaSheet.Range("P1").Copy()
$wobaSheet.Range("Q1:R1").PasteSpecial(-4122)

$wobaSheet.Range("Q2").Value = 91255
$wobaSheet.Range("R2").Value = 91663
$wobaSheet.Range("P2").Copy()
$wobaSheet.Range("Q2:R2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) "Park Factors" sheet: insert League column (D) after Abbreviation (C)
# ---------------------------------------------------------------------------
$pfSheet = $wb.Worksheets.Item("Park Factors")

$pfSheet.Range("D1").EntireColumn.Insert()
$pfSheet.Range("D1").Value = "League"

$leagueByRow = @{
    2 = "AL";  3 = "AL";  4 = "AL";  5 = "AL";  6 = "AL";
    7 = "AL";  8 = "AL";  9 = "AL";  10 = "AL"; 11 = "AL";
    12 = "AL"; 13 = "AL"; 14 = "AL"; 15 = "AL";
    16 = "NL"; 17 = "NL"; 18 = "NL"; 19 = "NL"; 20 = "NL"; 21 = "NL";
    22 = "AL";
    23 = "NL"; 24 = "NL"; 25 = "NL"; 26 = "NL"; 27 = "NL";
    28 = "NL"; 29 = "NL"; 30 = "NL"; 31 = "NL"
}

foreach ($row in 2..31) {
    $pfSheet.Cells.Item($row, 4).Value = $leagueByRow[$row]
}
